# Daily attendance processing - 2026-01-04 15:32:56
# Swap the "Recorded By" (column G) credit order for every session row that
# was recorded by both the System and dnasr281@gmail.com: the cell text
# "dnasr281@gmail.com, System" becomes "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "dnasr281@gmail.com, System"
$newText = "System, dnasr281@gmail.com"

$used = $ws.UsedRange
$lastRow = $used.Rows.Count
$colG = 7

$changed = 0
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colG)
    if ($cell.Text -eq $oldText) {
        $cell.Value = $newText
        $changed++
    }
}

Write-Output "Updated $changed cells in column G"
